$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-25 Thursday" "2024-07-26 Friday"

Replace-Text "498×9=" "691×2="
Replace-Text "751×4=" "764×3="
Replace-Text "209×4=" "835×3="
Replace-Text "495×3=" "901×6="
Replace-Text "356×6=" "300×3="

Replace-Text "380×2=" "778×6="
Replace-Text "865×7=" "897×3="
Replace-Text "736×3=" "322×9="
Replace-Text "719×9=" "383×2="
Replace-Text "948×4=" "515×2="

Replace-Text "647×4=" "121×3="
Replace-Text "320×3=" "975×5="
Replace-Text "591×8=" "151×6="
Replace-Text "863×4=" "251×3="
Replace-Text "535×8=" "251×2="

Replace-Text "561×7=" "831×2="
Replace-Text "489×8=" "370×9="
Replace-Text "469×2=" "802×4="
Replace-Text "457×4=" "421×3="
Replace-Text "670×2=" "822×4="

Replace-Text "135×8=" "177×2="
Replace-Text "106×6=" "868×7="
Replace-Text "780×2=" "183×7="
Replace-Text "116×5=" "763×6="
Replace-Text "863×2=" "910×9="
